$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.025136354050514
$ws.Cells.Item(2, 4).Value = 1.035935189204386
$ws.Cells.Item(2, 5).Value = 0.992614727750844
$ws.Cells.Item(2, 6).Value = 1.044369296142467
$ws.Cells.Item(2, 9).Value = 1.034971433796088
$ws.Cells.Item(2, 10).Value = 1.030307062009083
$ws.Cells.Item(2, 11).Value = 1.038730370075007
$ws.Cells.Item(2, 12).Value = 0.9955398523335997
$ws.Cells.Item(2, 13).Value = 1.047140539790082
$ws.Cells.Item(2, 14).Value = 1.014160805081473
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025968897501893
$ws.Cells.Item(3, 4).Value = 1.036587468207921
$ws.Cells.Item(3, 5).Value = 0.9936372048519299
$ws.Cells.Item(3, 6).Value = 1.045153782415341
$ws.Cells.Item(3, 9).Value = 1.035126499410188
$ws.Cells.Item(3, 10).Value = 1.030779414340665
$ws.Cells.Item(3, 11).Value = 1.039192564359862
$ws.Cells.Item(3, 12).Value = 0.9963617723202687
$ws.Cells.Item(3, 13).Value = 1.047736319840912
$ws.Cells.Item(3, 14).Value = 1.014316375590159
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.02650787317501
$ws.Cells.Item(4, 4).Value = 1.037009384041707
$ws.Cells.Item(4, 5).Value = 0.9942998659930998
$ws.Cells.Item(4, 6).Value = 1.045661579908158
$ws.Cells.Item(4, 9).Value = 1.035225096666726
$ws.Cells.Item(4, 10).Value = 1.031084676456742
$ws.Cells.Item(4, 11).Value = 1.039490798788019
$ws.Cells.Item(4, 12).Value = 0.9968940712668347
$ws.Cells.Item(4, 13).Value = 1.048121354518763
$ws.Cells.Item(4, 14).Value = 1.014416907200537
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.02673452032548
$ws.Cells.Item(5, 4).Value = 1.037186719161337
$ws.Cells.Item(5, 5).Value = 0.994578699834602
$ws.Cells.Item(5, 6).Value = 1.045875099696529
$ws.Cells.Item(5, 9).Value = 1.035266129350521
$ws.Cells.Item(5, 10).Value = 1.031212916110085
$ws.Cells.Item(5, 11).Value = 1.039615974926661
$ws.Cells.Item(5, 12).Value = 0.9971179600053012
$ws.Cells.Item(5, 13).Value = 1.048283107759554
$ws.Cells.Item(5, 14).Value = 1.014459138451795
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.026772578957826
$ws.Cells.Item(6, 4).Value = 1.037216492208928
$ws.Cells.Item(6, 5).Value = 0.994625531979634
$ws.Cells.Item(6, 6).Value = 1.045910952998343
$ws.Cells.Item(6, 9).Value = 1.035272994403753
$ws.Cells.Item(6, 10).Value = 1.031234442649829
$ws.Cells.Item(6, 11).Value = 1.039636980686582
$ws.Cells.Item(6, 12).Value = 0.9971555583673455
$ws.Cells.Item(6, 13).Value = 1.048310260028655
$ws.Cells.Item(6, 14).Value = 1.014466227362071
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.026510901404493
$ws.Cells.Item(7, 4).Value = 1.03701175375506
$ws.Cells.Item(7, 5).Value = 0.994303590798249
$ws.Cells.Item(7, 6).Value = 1.045664432808545
$ws.Cells.Item(7, 9).Value = 1.035225646589497
$ws.Cells.Item(7, 10).Value = 1.031086390366434
$ws.Cells.Item(7, 11).Value = 1.039492472191732
$ws.Cells.Item(7, 12).Value = 0.9968970624462089
$ws.Cells.Item(7, 13).Value = 1.048123516328778
$ws.Cells.Item(7, 14).Value = 1.014417471623785
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.025417660352534
$ws.Cells.Item(8, 4).Value = 1.036155660380381
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.044634378011409
$ws.Cells.Item(8, 9).Value = 1.035024198797315
$ws.Cells.Item(8, 10).Value = 1.030466773932786
$ws.Cells.Item(8, 11).Value = 1.038886743075496
$ws.Cells.Item(8, 12).Value = 0.9958175282591056
$ws.Cells.Item(8, 13).Value = 1.04734198377537
$ws.Cells.Item(8, 14).Value = 1.014213408132371
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.023493321702355
$ws.Cells.Item(9, 4).Value = 1.034646028820056
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.042820772697346
$ws.Cells.Item(9, 9).Value = 1.034655931265813
$ws.Cells.Item(9, 10).Value = 1.029372063236914
$ws.Cells.Item(9, 11).Value = 1.037813029761129
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.045961259606953
$ws.Cells.Item(9, 14).Value = 1.013852822650806
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.022211924876554
$ws.Cells.Item(10, 4).Value = 1.033638986145175
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.041612814394181
$ws.Cells.Item(10, 9).Value = 1.034401536851135
$ws.Cells.Item(10, 10).Value = 1.028640400364148
$ws.Cells.Item(10, 11).Value = 1.037093041607441
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.045038478183037
$ws.Cells.Item(10, 14).Value = 1.01361178493328
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.02165743917432
$ws.Cells.Item(11, 4).Value = 1.03320280077286
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.041090042881018
$ws.Cells.Item(11, 9).Value = 1.034289285457468
$ws.Cells.Item(11, 10).Value = 1.028323157141312
$ws.Cells.Item(11, 11).Value = 1.036780305081029
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.044638378802639
$ws.Cells.Item(11, 14).Value = 1.013507264527234
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.021451535181898
$ws.Cells.Item(12, 4).Value = 1.033040764496031
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.040895906612078
$ws.Cells.Item(12, 9).Value = 1.034247276107394
$ws.Cells.Item(12, 10).Value = 1.028205255714649
$ws.Cells.Item(12, 11).Value = 1.036663995428975
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.044489686228627
$ws.Cells.Item(12, 14).Value = 1.013468418943591
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.021495699698637
$ws.Cells.Item(13, 4).Value = 1.033075522595303
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.040937547463221
$ws.Cells.Item(13, 9).Value = 1.034256301459399
$ws.Cells.Item(13, 10).Value = 1.02823054881035
$ws.Cells.Item(13, 11).Value = 1.036688950800778
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.04452158478748
$ws.Cells.Item(13, 14).Value = 1.013476752444484
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.021640417921256
$ws.Cells.Item(14, 4).Value = 1.033189407152439
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.041073994600717
$ws.Cells.Item(14, 9).Value = 1.034285819355818
$ws.Cells.Item(14, 10).Value = 1.028313412654358
$ws.Cells.Item(14, 11).Value = 1.036770693846645
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.044626089406871
$ws.Cells.Item(14, 14).Value = 1.013504053986658
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.02172959114962
$ws.Cells.Item(15, 4).Value = 1.033259572907462
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.041158070119208
$ws.Cells.Item(15, 9).Value = 1.034303964699058
$ws.Cells.Item(15, 10).Value = 1.028364459468574
$ws.Cells.Item(15, 11).Value = 1.036821039211988
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.044690467884364
$ws.Cells.Item(15, 14).Value = 1.013520872458029
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.022248732135704
$ws.Cells.Item(16, 4).Value = 1.033667931763344
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.041647515131741
$ws.Cells.Item(16, 9).Value = 1.034408942499636
$ws.Cells.Item(16, 10).Value = 1.028661445828383
$ws.Cells.Item(16, 11).Value = 1.037113776437036
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.045065020455337
$ws.Cells.Item(16, 14).Value = 1.013618718493742
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.022574475004385
$ws.Cells.Item(17, 4).Value = 1.033924051161647
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.041954608022031
$ws.Cells.Item(17, 9).Value = 1.034474231497561
$ws.Cells.Item(17, 10).Value = 1.028847623679708
$ws.Cells.Item(17, 11).Value = 1.037297142200624
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.045299826920324
$ws.Cells.Item(17, 14).Value = 1.013680054993798
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.022764510676762
$ws.Cells.Item(18, 4).Value = 1.034073428681969
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.042133757228873
$ws.Cells.Item(18, 9).Value = 1.034512111084621
$ws.Cells.Item(18, 10).Value = 1.02895617665944
$ws.Cells.Item(18, 11).Value = 1.037404001998417
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.045436734416212
$ws.Cells.Item(18, 14).Value = 1.013715817086322
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.022829313884652
$ws.Cells.Item(19, 4).Value = 1.034124360365215
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.042194847020206
$ws.Cells.Item(19, 9).Value = 1.034524992691504
$ws.Cells.Item(19, 10).Value = 1.028993183359979
$ws.Cells.Item(19, 11).Value = 1.037440422382291
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.045483407610353
$ws.Cells.Item(19, 14).Value = 1.013728008572275
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.022539522208494
$ws.Cells.Item(20, 4).Value = 1.033896573263403
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.041921657042211
$ws.Cells.Item(20, 9).Value = 1.034467247532967
$ws.Cells.Item(20, 10).Value = 1.028827652832959
$ws.Cells.Item(20, 11).Value = 1.037277478538544
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.045274639687025
$ws.Cells.Item(20, 14).Value = 1.013673475659532
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.021597800435169
$ws.Cells.Item(21, 4).Value = 1.033155871460445
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.041033813079994
$ws.Cells.Item(21, 9).Value = 1.034277135737212
$ws.Cells.Item(21, 10).Value = 1.028289013051853
$ws.Cells.Item(21, 11).Value = 1.036746626559162
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.044595317537252
$ws.Cells.Item(21, 14).Value = 1.013496014968416
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.021006031268378
$ws.Cells.Item(22, 4).Value = 1.032690061851675
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.040475847868029
$ws.Cells.Item(22, 9).Value = 1.034155787551968
$ws.Cells.Item(22, 10).Value = 1.027949983955377
$ws.Cells.Item(22, 11).Value = 1.036412018399598
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.044167751120431
$ws.Cells.Item(22, 14).Value = 1.013384311008266
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.021319707657002
$ws.Cells.Item(23, 4).Value = 1.032937005308963
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.040771610734505
$ws.Cells.Item(23, 9).Value = 1.03422028850035
$ws.Cells.Item(23, 10).Value = 1.028129743885188
$ws.Cells.Item(23, 11).Value = 1.036589479719429
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.04439445424847
$ws.Cells.Item(23, 14).Value = 1.013443539329681
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.022555315763296
$ws.Cells.Item(24, 4).Value = 1.033908989381958
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.041936546088053
$ws.Cells.Item(24, 9).Value = 1.034470403911352
$ws.Cells.Item(24, 10).Value = 1.02883667692705
$ws.Cells.Item(24, 11).Value = 1.037286363992698
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.045286020871777
$ws.Cells.Item(24, 14).Value = 1.013676448622175
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.023990552281296
$ws.Cells.Item(25, 4).Value = 1.035036422228206
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.043289445546633
$ws.Cells.Item(25, 9).Value = 1.034752706772204
$ws.Cells.Item(25, 10).Value = 1.029655404541185
$ws.Cells.Item(25, 11).Value = 1.038091353154891
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.046318622167365
$ws.Cells.Item(25, 14).Value = 1.013946158546376
